$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Is my Order ID 1588 refundable?"
$ws.Range("C2").Value = "{'orders': [{'Product ID': 21448, 'ProductName': 'panasonic dmc fz1000eb lumix bridge camera 25 400mm leica dc lens 20.1mp', 'Category': 'Digital Cameras', 'Category ID': 2617, 'OrderID': 1588, 'CustomerID': 868, 'OrderStatus': 'Pending', 'ReturnEligible': False, 'ShippingDate': '2024-10-08 12:24:37.765490'}]}"

$ws.Range("B3").Value = "Can I exchange my Order ID 2443?"
$ws.Range("C3").Value = "{'orders': [{'Product ID': 17894, 'ProductName': 'intel bx80532ke3066e processor 3.06 ghz 1 mb l2', 'Category': 'CPUs', 'Category ID': 2615, 'OrderID': 2443, 'CustomerID': 2440, 'OrderStatus': 'Pending', 'ReturnEligible': False, 'ShippingDate': '2024-10-14 12:24:37.765490'}]}"

$ws.Range("B4").Value = "What's the return policy for Order ID 2212?"
$ws.Range("C4").Value = "{'orders': [{'Product ID': 13535, 'ProductName': '10 hd portable television august da100d', 'Category': 'TVs', 'Category ID': 2614, 'OrderID': 2212, 'CustomerID': 1191, 'OrderStatus': 'Delivered', 'ReturnEligible': True, 'ShippingDate': '2024-09-29 12:24:37.765490'}]}"

$ws.Range("B5").Value = "My Order ID 1425 arrived damaged, can I return it?"
$ws.Range("C5").Value = "{'orders': [{'Product ID': 3714, 'ProductName': 'samsung galaxy s ii white mobile phone', 'Category': 'Mobile Phones', 'Category ID': 2612, 'OrderID': 1425, 'CustomerID': 2200, 'OrderStatus': 'Delivered', 'ReturnEligible': False, 'ShippingDate': '2024-10-09 12:24:37.765490'}]}"

$ws.Range("B6").Value = "I received the wrong item in Order ID 1678, how do I return it?"
$ws.Range("C6").Value = "{'orders': [{'Product ID': 37900, 'ProductName': 'lg gsl761wbxv american fridge freezer in black ice water a rated', 'Category': 'Fridge Freezers', 'Category ID': 2622, 'OrderID': 1678, 'CustomerID': 947, 'OrderStatus': 'Pending', 'ReturnEligible': False, 'ShippingDate': '2024-10-03 12:24:37.765490'}]}"

$ws.Range("B7").Value = "Can I return part of my Order ID 1189?"
$ws.Range("C7").Value = "{'orders': [{'Product ID': 1140, 'ProductName': 'sim free nokia 3.1 16gb mobile phone black/silver', 'Category': 'Mobile Phones', 'Category ID': 2612, 'OrderID': 1189, 'CustomerID': 577, 'OrderStatus': 'Pending', 'ReturnEligible': False, 'ShippingDate': '2024-10-23 12:24:37.765490'}]}"

$ws.Range("B8").Value = "Who pays for return shipping for Order ID 3279?"
$ws.Range("C8").Value = "{'orders': [{'Product ID': 12483, 'ProductName': 'lg electronics 50uk6950 fh fernseher ultra hd 4k', 'Category': 'TVs', 'Category ID': 2614, 'OrderID': 3279, 'CustomerID': 873, 'OrderStatus': 'Delivered', 'ReturnEligible': True, 'ShippingDate': '2024-10-09 12:24:37.765490'}]}"

$ws.Range("B9").Value = "Is there a restocking fee for returning Order ID 3124?"
$ws.Range("C9").Value = "{'orders': [{'Product ID': 22525, 'ProductName': 'canon eos 1300d slr camera inc ef s 18 55mm f/3.5 5.6 is ii lens', 'Category': 'Digital Cameras', 'Category ID': 2617, 'OrderID': 3124, 'CustomerID': 803, 'OrderStatus': 'Pending', 'ReturnEligible': False, 'ShippingDate': '2024-10-02 12:24:37.765490'}]}"

$ws.Range("B10").Value = "Can I return Order ID 1530 from another country?"
$ws.Range("C10").Value = "{'orders': [{'Product ID': 2352, 'ProductName': 'wiko jerry schwarz grau', 'Category': 'Mobile Phones', 'Category ID': 2612, 'OrderID': 1530, 'CustomerID': 812, 'OrderStatus': 'Delivered', 'ReturnEligible': True, 'ShippingDate': '2024-10-15 12:24:37.765490'}]}"

$ws.Range("B11").Value = "Is my order 2594 refundable?"
$ws.Range("C11").Value = "{'orders': [{'Product ID': 46656, 'ProductName': 'amica einbau k hlschrank schleppt rtechnik eks 16171', 'Category': 'Fridges', 'Category ID': 2623, 'OrderID': 2286, 'CustomerID': 46, 'OrderStatus': 'Shipped', 'ReturnEligible': False, 'ShippingDate': '2024-10-07 12:24:37.765490'}]}"
